# Generate Report for Handback
# Updates the localization-status report: marks the handback rows as
# "in sync with en-US", refreshes the handback timestamps, clears the
# (now resolved) error-detail messages, and widens a couple of report
# columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "Handed back: in sync with en-US"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Columns("E").ColumnWidth = 29.167
$ws.Columns("F").ColumnWidth = 29.167

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K2").Value = "2016-08-13 22:47:38"
$ws.Range("K3").Value = "2016-08-13 22:47:38"
$ws.Range("P2").Value = ""
$ws.Range("P3").Value = ""
$ws.Columns("C").ColumnWidth = 29.167
$ws.Columns("P").ColumnWidth = 12.8335

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K2").Value = "2016-08-13 22:47:48"
$ws.Range("K3").Value = "2016-08-13 22:47:48"
$ws.Range("P2").Value = ""
$ws.Range("P3").Value = ""
$ws.Columns("C").ColumnWidth = 29.167
$ws.Columns("P").ColumnWidth = 12.8335
